$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# Row 24: mark as "Hecho", set hours estimate and day-10 consumption
$ws.Range("F24").Value = "Hecho"
$ws.Range("G24").Value = 1
$ws.Range("AI24").Value = 1

# Row 26: mark as "Hecho"
$ws.Range("F26").Value = "Hecho"

# Row 27: mark as "Hecho", set day-11 consumption
$ws.Range("F27").Value = "Hecho"
$ws.Range("AL27").Value = 2

# Row 28: mark as "Hecho", set hours estimate and day-11 consumption
$ws.Range("F28").Value = "Hecho"
$ws.Range("G28").Value = 1
$ws.Range("AL28").Value = 1

# Row 29: mark as "Hecho", set day-11 consumption
$ws.Range("F29").Value = "Hecho"
$ws.Range("AL29").Value = 1

# Row 30: mark as "En proceso", move consumption from day-10 to day-11
$ws.Range("F30").Value = "En proceso"
$ws.Range("AI30").ClearContents()
$ws.Range("AL30").Value = 1

# Update the active selection/scroll position to reflect the edited area
$ws.Activate()
$ws.Range("AI29").Select()
